# L-D5019-01-05_EVE.xlsx: drop C8 / C10 from the 22p/0805 capacitor group
# (SPI can be clocked a little faster without them), and adjust the
# related quantity / row-height / selection bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Designator list for the 22p / 0805 capacitor row (row 3): remove C8, C10.
$ws.Range("B3").Value = "C2, C3, C4, C5"

# Quantity for that row drops from 6 to 4 designators.
$ws.Range("A3").Value = 4

# Row 3 height was nudged slightly when the sheet was re-saved.
$ws.Rows(3).RowHeight = 13.4

# Cursor/selection ends up parked on B4 instead of A23.
[void]$ws.Range("B4").Select()
